$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New query strings for each tab (column B) and the shared StatQuery (column C) ---

$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (samp:sample)-->(c) 
  MATCH (f:file)-[*]->(c)
   WHERE f.file_type IN ["RNA Sequence File"] 
OPTIONAL MATCH (co:cohort)<-[*]-(c)
  WITH DISTINCT c, s, demo, diag, co,demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
  coalesce(CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '') AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
order by c.case_id asc
limit 100
'@

$statQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(f)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE  f.file_type IN ["RNA Sequence File"]   
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$samplesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
MATCH (f:file)-[*]->(c)
WHERE f.file_type IN  ["RNA Sequence File"]    
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed,
        coalesce(diag.disease_term,'') AS Diagnosis, 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
order by samp.sample_id asc
limit 100
'@

$filesQuery = @'
MATCH (f:file)-[*]->(c:case)
MATCH (f)-->(parent)
OPTIONAL MATCH (f)-->(samp:sample)
OPTIONAL MATCH (p:program)<--(s:study)<--(c)
OPTIONAL MATCH (c)<--(demo:demographic)
OPTIONAL MATCH (c)<--(diag:diagnosis)
WITH
    f, c, parent, samp, p, s, demo, diag
WHERE f.file_type IN ["RNA Sequence File"]  
WITH
    DISTINCT f, c, parent, samp, p, s, demo, diag,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    2 as precision
WITH
    DISTINCT f, c, parent, samp, p, s, demo, diag,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    DISTINCT f, c, parent, samp, p, s, demo, diag, unit,
    round(factor * value)/factor AS size
RETURN
    coalesce(f.file_name, '') AS `File Name`,
    coalesce(f.file_format, '') AS `Format`,
    coalesce(f.file_type, '') AS `File Type`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    coalesce(labels(parent)[0], '') AS `Association`,
    coalesce(f.file_description, '') AS `Description`,
    coalesce(samp.sample_id, '') AS `Sample ID`,
    coalesce(c.case_id, '') AS `Case ID`,
    coalesce(demo.breed,'') AS Breed ,
    coalesce(diag.disease_term,'') AS Diagnosis
    ORDER BY f.file_name asc
     limit 100
'@

$studyFilesQuery = @'
MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(diag:diagnosis)
MATCH (samp:sample)-->(c)
MATCH (c)<--(demo:demographic)
WHERE f.file_type IN ["RNA Sequence File"] 
WITH
        DISTINCT f, c, demo, diag, s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
  order by 'File Name' asc
  limit 100
'@

# --- Row 2: CasesTab ---
$ws.Range("B2").Value = $casesQuery
$ws.Range("C2").Value = $statQuery

# --- Row 3: SamplesTab ---
$ws.Range("B3").Value = $samplesQuery
$ws.Range("C3").Value = $statQuery

# --- Row 4: FilesTab ---
$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $statQuery

# --- Row 5 (new): StudyFilesTab ---
$ws.Range("A5").Value = "StudyFilesTab"
$ws.Range("B5").Value = $studyFilesQuery
$ws.Range("C5").Value = $statQuery
$ws.Range("D5").Value = $ws.Range("D2").Value
$ws.Range("E5").Value = $ws.Range("E2").Value

# --- Formatting: bump font size to 16 across the used range, and wrap the query columns ---
$ws.Range("A1:E5").Font.Size = 16
$ws.Range("B2:C5").WrapText = $true

# --- Row heights: every data row is now fully expanded (Excel's max row height) ---
$ws.Rows(2).RowHeight = 409.5
$ws.Rows(3).RowHeight = 409.5
$ws.Rows(4).RowHeight = 409.5
$ws.Rows(5).RowHeight = 409.5

# --- Selection moves to B10 as in the authored workbook ---
$ws.Range("B10").Select()

Write-Host "edit applied"
